$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 9885.1875
$ws.Range("I11").Value = 9885.1875
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 9885.1875
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -9745.1875

$ws.Range("H100").Value = 3745.3333
$ws.Range("I100").Value = 3890
$ws.Range("J100").Value = 3166.6667
$ws.Range("K100").Value = 3890
$ws.Range("L100").Value = 3166.6667
$ws.Range("M100").Value = -3349

$ws.Range("H116").Value = 4886
$ws.Range("I116").Value = 4622.375
$ws.Range("J116").Value = 6995
$ws.Range("K116").Value = 4622.375
$ws.Range("L116").Value = 6995
$ws.Range("M116").Value = -1180.375
$ws.Range("N116").Value = -13879

$ws.Range("H137").Value = 5078.6743
$ws.Range("I137").Value = 2626.16
$ws.Range("J137").Value = 8484.944
$ws.Range("K137").Value = 7878.48
$ws.Range("L137").Value = 25454.832
$ws.Range("M137").Value = -5328.48
$ws.Range("N137").Value = -30554.832

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15262.797
$ws.Range("I32").Value = 11481.882
$ws.Range("J32").Value = 59877.6
$ws.Range("K32").Value = 11481.882
$ws.Range("L32").Value = 59877.6
$ws.Range("M32").Value = -11194.882
$ws.Range("N32").Value = -60451.6

$ws.Range("H61").Value = 4560.8965
$ws.Range("I61").Value = 3216.8667
$ws.Range("J61").Value = 6000.9287
$ws.Range("K61").Value = 3216.8667
$ws.Range("L61").Value = 6000.9287
$ws.Range("M61").Value = -3004.8667
$ws.Range("N61").Value = -6424.9287

$ws.Range("H102").Value = 1062.0834
$ws.Range("I102").Value = 1062.0834
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1062.0834
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 559.9166

$ws.Range("H110").Value = 3735.75
$ws.Range("I110").Value = 3555.1428
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 3555.1428
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -1510.1428

$ws.Range("H136").Value = 4560.8965
$ws.Range("I136").Value = 3216.8667
$ws.Range("J136").Value = 6000.9287
$ws.Range("K136").Value = 9650.6001
$ws.Range("L136").Value = 18002.7861
$ws.Range("M136").Value = -7100.6001
$ws.Range("N136").Value = -23102.7861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2004.8462
$ws.Range("I20").Value = 1856.4
$ws.Range("J20").Value = 2499.6667
$ws.Range("K20").Value = 1856.4
$ws.Range("L20").Value = 2499.6667
$ws.Range("M20").Value = -1609.4
$ws.Range("N20").Value = -2993.6667

$ws.Range("H42").Value = 187142.72
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 187142.72
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 187142.72
$ws.Range("N42").Value = -187798.72

$ws.Range("H99").Value = 4092.3333
$ws.Range("I99").Value = 3888.5
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 3888.5
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -2390.5
$ws.Range("N99").Value = -7496

$ws.Range("H105").Value = 3695.6553
$ws.Range("I105").Value = 3363.9333
$ws.Range("J105").Value = 4051.0715
$ws.Range("K105").Value = 3363.9333
$ws.Range("L105").Value = 4051.0715
$ws.Range("M105").Value = -1616.9333

$ws.Range("H107").Value = 3796.6667
$ws.Range("I107").Value = 3796.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3796.6667
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1876.6667

$ws.Range("H134").Value = 6999.1665
$ws.Range("I134").Value = 5997.5
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 17992.5
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = -15457.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9018.857
$ws.Range("I62").Value = 8971.362999999999
$ws.Range("J62").Value = 9193
$ws.Range("K62").Value = 8971.362999999999
$ws.Range("L62").Value = 9193
$ws.Range("M62").Value = -8347.362999999999

$ws.Range("H65").Value = 9018.857
$ws.Range("I65").Value = 8971.362999999999
$ws.Range("J65").Value = 9193
$ws.Range("K65").Value = 44856.815
$ws.Range("L65").Value = 45965
$ws.Range("M65").Value = -41736.815

$ws.Range("H92").Value = 49000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 49000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 49000
$ws.Range("N92").Value = -53992

$ws.Range("H141").Value = 535817.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 535817.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 535817.25
$ws.Range("N141").Value = -546177.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 269
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 269
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 807
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -1031

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H12").Value = 158
$ws.Range("I12").Value = 46.8
$ws.Range("J12").Value = 227.5
$ws.Range("K12").Value = 140.4
$ws.Range("L12").Value = 682.5
$ws.Range("M12").Value = 32.60000000000002
$ws.Range("N12").Value = -1028.5

$ws.Range("H68").Value = 998
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 998
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2994
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4616

$ws.Range("H71").Value = 998
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 998
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 8982
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -17094

$ws.Range("H81").Value = 110012
$ws.Range("I81").Value = 110013
$ws.Range("J81").Value = 110011.5
$ws.Range("K81").Value = 330039
$ws.Range("L81").Value = 330034.5
$ws.Range("M81").Value = -328916
$ws.Range("N81").Value = -332280.5

$ws.Range("H84").Value = 110012
$ws.Range("I84").Value = 110013
$ws.Range("J84").Value = 110011.5
$ws.Range("K84").Value = 990117
$ws.Range("L84").Value = 990103.5
$ws.Range("M84").Value = -984501
$ws.Range("N84").Value = -1001335.5

$ws.Range("H92").Value = 499
$ws.Range("I92").Value = 499
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1497
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -249

$ws.Range("H97").Value = 811.3333
$ws.Range("I97").Value = 811.3333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2433.9999
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1937.9999
$ws.Range("N97").ClearContents()

$ws.Range("H113").Value = 1216.65
$ws.Range("I113").Value = 999.25
$ws.Range("J113").Value = 1271
$ws.Range("K113").Value = 2997.75
$ws.Range("L113").Value = 3813
$ws.Range("M113").Value = -827.75
$ws.Range("N113").Value = -8153

$ws.Range("H122").Value = 1402.56
$ws.Range("I122").Value = 1298.3636
$ws.Range("J122").Value = 1484.4286
$ws.Range("K122").Value = 11685.2724
$ws.Range("L122").Value = 13359.8574
$ws.Range("M122").Value = -9235.2724
$ws.Range("N122").Value = -18259.8574

$ws.Range("H132").Value = 2560.1562
$ws.Range("I132").Value = 1449.1111
$ws.Range("J132").Value = 2994.913
$ws.Range("K132").Value = 13041.9999
$ws.Range("L132").Value = 26954.217
$ws.Range("M132").Value = -10511.9999
$ws.Range("N132").Value = -32014.217

$ws.Range("H135").Value = 269
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 269
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 2421
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -7491

$ws.Range("H139").Value = 8011.3516
$ws.Range("I139").Value = 5932.385
$ws.Range("J139").Value = 9137.458000000001
$ws.Range("K139").Value = 17797.155
$ws.Range("L139").Value = 27412.374
$ws.Range("M139").Value = -12657.155
$ws.Range("N139").Value = -37692.374

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6513.8335
$ws.Range("I132").Value = 6606.5454
$ws.Range("J132").Value = 5494
$ws.Range("K132").Value = 19819.6362
$ws.Range("L132").Value = 16482
$ws.Range("M132").Value = -17289.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 40000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 40000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352

$ws.Range("H136").Value = 3764.375
$ws.Range("I136").Value = 3719.8333
$ws.Range("J136").Value = 3898
$ws.Range("K136").Value = 11159.4999
$ws.Range("L136").Value = 11694
$ws.Range("M136").Value = -8609.499899999999
$ws.Range("N136").Value = -16794

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 39994
$ws.Range("I54").Value = 39994
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 39994
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -39474

$ws.Range("H132").Value = 1910.7142
$ws.Range("I132").Value = 1673.0769
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5019.2307
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2489.2307
$ws.Range("N132").Value = -20060

$ws.Range("H136").Value = 2489.239
$ws.Range("I136").Value = 2009.8667
$ws.Range("J136").Value = 3388.0625
$ws.Range("K136").Value = 6029.6001
$ws.Range("L136").Value = 10164.1875
$ws.Range("M136").Value = -3479.6001
